$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RecoveryBoardBOM")

# Fill in the Digi-Key part info for the two "10k" resistor rows (R1,R4,R16 and
# R7,R10,... group), matching the existing 100k resistor row (DPN/MFR/MPN).
$ws.Range("H43").Value = "311-10KGRCT-ND"
$ws.Range("I43").Value = "Yageo"
$ws.Range("J43").Value = "RC0603JR-0710KL"

$ws.Range("H49").Value = "311-10KGRCT-ND"
$ws.Range("I49").Value = "Yageo"
$ws.Range("J49").Value = "RC0603JR-0710KL"

# Update the frozen-pane scroll position and current selection to reflect
# where the user ended up after editing (bottom of the sheet).
$ws.Activate()
$ws.Range("A28").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("G49:J49").Select()
